# Locale.xlsx update: "Tooltip de botin ordenado y filtrado"
# Inserts/reorders localization key-value pairs (rows 20-79) to match
# the new sharedStrings ordering, and adds new item-rarity tooltip strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(20, 1).Value = 'Press key to show'
$ws.Cells.Item(20, 2).Value = 'Pulse la tecla para mostrar'
$ws.Cells.Item(21, 1).Value = 'Accept'
$ws.Cells.Item(21, 2).Value = 'Aceptar'
$ws.Cells.Item(22, 1).Value = 'Cancel'
$ws.Cells.Item(22, 2).Value = 'Cancelar'
$ws.Cells.Item(23, 1).Value = 'Module |cffffcc00%s|r loaded'
$ws.Cells.Item(23, 2).Value = 'Módulo |cffffcc00%s|r cargado'
$ws.Cells.Item(24, 1).Value = 'LogBook %s initialized'
$ws.Cells.Item(24, 2).Value = 'LogBook %s inicializado'
$ws.Cells.Item(25, 1).Value = 'Left Click'
$ws.Cells.Item(25, 2).Value = 'Clic Izquierdo'
$ws.Cells.Item(26, 1).Value = 'Open main window'
$ws.Cells.Item(26, 2).Value = 'Abrir la ventana principal'
$ws.Cells.Item(27, 1).Value = 'Right Click'
$ws.Cells.Item(27, 2).Value = 'Clic Derecho'
$ws.Cells.Item(28, 1).Value = 'Open settings window'
$ws.Cells.Item(28, 2).Value = 'Abrir la ventana de configuración'
$ws.Cells.Item(29, 1).Value = 'Advanced'
$ws.Cells.Item(29, 2).Value = 'Avanzado'
$ws.Cells.Item(30, 1).Value = 'Advanced settings'
$ws.Cells.Item(30, 2).Value = 'Ajustes avanzados'
$ws.Cells.Item(31, 1).Value = 'Debug'
$ws.Cells.Item(31, 2).Value = 'Depurar'
$ws.Cells.Item(32, 1).Value = 'Enable debug'
$ws.Cells.Item(32, 2).Value = 'Habilitar depuración'
$ws.Cells.Item(33, 1).Value = 'General'
$ws.Cells.Item(33, 2).Value = 'General'
$ws.Cells.Item(34, 1).Value = 'General settings'
$ws.Cells.Item(34, 2).Value = 'Configuración general'
$ws.Cells.Item(35, 1).Value = 'Log|cff57b6ffBook|r available commands'
$ws.Cells.Item(35, 2).Value = 'Comandos disponibles de Log|cff57b6ffBook|r'
$ws.Cells.Item(36, 1).Value = 'Open critics window'
$ws.Cells.Item(36, 2).Value = 'Abrir ventana de críticos'
$ws.Cells.Item(37, 1).Value = 'Open loot window'
$ws.Cells.Item(37, 2).Value = 'Abrir ventana de botín'
$ws.Cells.Item(38, 1).Value = 'Open zones window'
$ws.Cells.Item(38, 2).Value = 'Abrir la ventana de zonas'
$ws.Cells.Item(39, 1).Value = 'Open fishing window'
$ws.Cells.Item(39, 2).Value = 'Abrir la ventana de pesca'
$ws.Cells.Item(40, 1).Value = 'Open mobs window'
$ws.Cells.Item(40, 2).Value = 'Abrir la ventana de mobs'
$ws.Cells.Item(41, 1).Value = 'Open enchanting window'
$ws.Cells.Item(41, 2).Value = 'Abrir la ventana de encantamiento'
$ws.Cells.Item(42, 1).Value = 'Main plugins'
$ws.Cells.Item(42, 2).Value = 'Plugins principales'
$ws.Cells.Item(43, 1).Value = 'Loot'
$ws.Cells.Item(43, 2).Value = 'Botín'
$ws.Cells.Item(44, 1).Value = 'Allows you to track loot and items crafted with trading skills.'
$ws.Cells.Item(44, 2).Value = 'Permite realizar un seguimiento del botín y los artículos elaborados con habilidades comerciales.'
$ws.Cells.Item(45, 1).Value = 'Fishing'
$ws.Cells.Item(45, 2).Value = 'Pesca'
$ws.Cells.Item(46, 1).Value = 'Allows you to track fish from pools and wreckages.'
$ws.Cells.Item(46, 2).Value = 'Permite realizar un seguimiento de peces desde pozas y restos de naufragios.'
$ws.Cells.Item(47, 1).Value = 'Critics'
$ws.Cells.Item(47, 2).Value = 'Críticos'
$ws.Cells.Item(48, 1).Value = 'Allows you to track hits or healing, both normal and critical.'
$ws.Cells.Item(48, 2).Value = 'Permite realizar un seguimiento de los golpes o curaciones, tanto normales como críticos.'
$ws.Cells.Item(49, 1).Value = 'Zones'
$ws.Cells.Item(49, 2).Value = 'Zonas'
$ws.Cells.Item(50, 1).Value = 'Allows you to track zones.'
$ws.Cells.Item(50, 2).Value = 'Permite realizar un seguimiento de las zonas.'
$ws.Cells.Item(51, 1).Value = 'Mobs'
$ws.Cells.Item(51, 2).Value = 'Mobs'
$ws.Cells.Item(52, 1).Value = 'Allows you to track mobs.'
$ws.Cells.Item(52, 2).Value = 'Permite realizar un seguimiento de las mobs.'
$ws.Cells.Item(53, 1).Value = 'Enchanting'
$ws.Cells.Item(53, 2).Value = 'Encantamiento'
$ws.Cells.Item(54, 1).Value = 'Allows you to track enchanting.'
$ws.Cells.Item(54, 2).Value = 'Permite realizar un seguimiento de encantamiento.'
$ws.Cells.Item(55, 1).Value = 'Settings'
$ws.Cells.Item(55, 2).Value = 'Configuración'
$ws.Cells.Item(56, 1).Value = 'Stats'
$ws.Cells.Item(56, 2).Value = 'Estadísticas'
$ws.Cells.Item(57, 1).Value = 'Tooltips'
$ws.Cells.Item(57, 2).Value = 'Mensaje emergente'
$ws.Cells.Item(58, 1).Value = 'Maintenance'
$ws.Cells.Item(58, 2).Value = 'Mantenimiento'
$ws.Cells.Item(59, 1).Value = 'Starting database auto update: %s'
$ws.Cells.Item(59, 2).Value = 'Iniciando la actualización automática de la base de datos: %s'
$ws.Cells.Item(60, 1).Value = 'Cancelling database auto update: %s'
$ws.Cells.Item(60, 2).Value = 'Cancelando la actualización automática de la base de datos: %s'
$ws.Cells.Item(61, 1).Value = '%s database update: %s'
$ws.Cells.Item(61, 2).Value = 'Actualización base de datos %s : %s'
$ws.Cells.Item(62, 1).Value = 'Done!'
$ws.Cells.Item(62, 2).Value = '¡Hecho!'
$ws.Cells.Item(63, 1).Value = 'All'
$ws.Cells.Item(63, 2).Value = 'Todos'
$ws.Cells.Item(64, 1).Value = 'Database'
$ws.Cells.Item(64, 2).Value = 'Base de datos'
$ws.Cells.Item(65, 1).Value = 'more'
$ws.Cells.Item(65, 2).Value = 'más'
$ws.Cells.Item(66, 1).Value = 'Normal'
$ws.Cells.Item(66, 2).Value = 'Normal'
$ws.Cells.Item(67, 1).Value = 'Rare'
$ws.Cells.Item(67, 2).Value = 'Raro'
$ws.Cells.Item(68, 1).Value = 'Elite'
$ws.Cells.Item(68, 2).Value = 'Élite'
$ws.Cells.Item(69, 1).Value = 'Rare elite'
$ws.Cells.Item(69, 2).Value = 'Raro élite'
$ws.Cells.Item(70, 1).Value = 'Boss'
$ws.Cells.Item(70, 2).Value = 'Jefe'
$ws.Cells.Item(71, 1).Value = 'World boss'
$ws.Cells.Item(71, 2).Value = 'Jefe de mundo'
$ws.Cells.Item(72, 1).Value = 'Poor'
$ws.Cells.Item(72, 2).Value = 'Pobre'
$ws.Cells.Item(73, 1).Value = 'Common'
$ws.Cells.Item(73, 2).Value = 'Común'
$ws.Cells.Item(74, 1).Value = 'Uncommon'
$ws.Cells.Item(74, 2).Value = 'Poco Común'
$ws.Cells.Item(75, 1).Value = 'Epic'
$ws.Cells.Item(75, 2).Value = 'Epico'
$ws.Cells.Item(76, 1).Value = 'Legendary'
$ws.Cells.Item(76, 2).Value = 'Legendario'
$ws.Cells.Item(77, 1).Value = 'Artifact'
$ws.Cells.Item(77, 2).Value = 'Artefacto'
$ws.Cells.Item(78, 1).Value = 'Heirloom'
$ws.Cells.Item(78, 2).Value = 'Reliquia'
$ws.Cells.Item(79, 1).Value = 'WoW Token'
$ws.Cells.Item(79, 2).Value = 'Ficha WoW'

# Update the active selection / view state to reflect where the author
# left the cursor after editing (row 61 area, column D selected).
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("D61").Select()
$ws.Columns("D:D").Select()
